$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.006.99"
$ws.Range("E2").Value = "  +2.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.252.88"
$ws.Range("E3").Value = "  +1.55%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.65"
$ws.Range("E5").Value = "  -0.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.67"
$ws.Range("E6").Value = "  +3.07%  "

$ws.Range("E7").Value = "  -0.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.552"
$ws.Range("E9").Value = "  -1.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.41"
$ws.Range("E10").Value = "  +1.90%  "

$ws.Range("E11").Value = "  +0.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.58"
$ws.Range("E12").Value = "  -0.43%  "

$ws.Range("E13").Value = "  -1.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.598.91"
$ws.Range("E14").Value = "  +1.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.856"
$ws.Range("E15").Value = "  -0.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.42"
$ws.Range("E16").Value = "  +0.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.259.47"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.908.33"
$ws.Range("E18").Value = "  +2.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.24"
$ws.Range("E19").Value = "  -5.43%  "

$ws.Range("E20").Value = "  +2.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.51"
$ws.Range("E21").Value = "  -0.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.70"
$ws.Range("E22").Value = "  +0.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.14"
$ws.Range("E23").Value = "  -1.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.31"
$ws.Range("E24").Value = "  -0.43%  "

$ws.Range("E25").Value = "  -2.29%  "

$ws.Range("E26").Value = "  +0.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.22"
$ws.Range("E27").Value = "  +2.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  +1.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.92"
$ws.Range("E29").Value = "  +6.59%  "

$ws.Range("E30").Value = "  -2.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "161.64"
$ws.Range("E31").Value = "  +5.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.18"
$ws.Range("E32").Value = "  -0.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0850"
$ws.Range("E33").Value = "  -1.69%  "

$ws.Range("E34").Value = "  +1.25%  "

$ws.Range("E35").Value = "  +10.99%  "

$ws.Range("E36").Value = "  +0.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.01"
$ws.Range("E37").Value = "  -6.22%  "

$ws.Range("E38").Value = "  -1.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.88"
$ws.Range("E39").Value = "  +23.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.71"
$ws.Range("E40").Value = "  +1.22%  "

$ws.Range("E41").Value = "  -4.22%  "

$ws.Range("E42").Value = "  -1.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.791.81"
$ws.Range("E44").Value = "  +3.82%  "

$ws.Range("E45").Value = "  -2.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "75.48"
$ws.Range("E46").Value = "  +1.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "82.37"
$ws.Range("E47").Value = "  -2.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.20"
$ws.Range("E48").Value = "  -1.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.75"
$ws.Range("E49").Value = "  +1.84%  "

$ws.Range("E50").Value = "  +8.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.21"
$ws.Range("E51").Value = "  +0.98%  "
